$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 updates
$ws.Range("B2").Value = $false
$ws.Range("C2").Value = 10000
$ws.Range("D2").Value = 19.32
$ws.Range("E2").Value = 19.170000000000002
$ws.Range("F2").Value = 0
$ws.Range("G2").Value = $false

# New row 3
$ws.Range("C3").Value = 0

# Update selection to match diff (active cell G2)
$ws.Range("G2").Select()
